$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update revised orshanky basket calculations (rows 2-29)
$ws.Range("C2").Value = 93.62039390430655
$ws.Range("D2").Value = 114.94969125460553
$ws.Range("E2").Value = 2976.0
$ws.Range("G2").Value = 396.5764465332031
$ws.Range("H2").Value = 7471.72998046875

$ws.Range("C3").Value = 98.85272879049342
$ws.Range("D3").Value = 121.37409532839253
$ws.Range("E3").Value = 2976.0
$ws.Range("G3").Value = 464.86279296875
$ws.Range("H3").Value = 9103.0576171875

$ws.Range("C4").Value = 10.393337275392266
$ws.Range("D4").Value = 12.761224837953685
$ws.Range("E4").Value = 2976.0
$ws.Range("G4").Value = 45.30234909057617
$ws.Range("H4").Value = 446.64288330078125

$ws.Range("C5").Value = 1.8349654508694526
$ws.Range("D5").Value = 2.253020942852061
$ws.Range("E5").Value = 2976.0
$ws.Range("G5").Value = 6.413599491119385
$ws.Range("H5").Value = 324.4350280761719

$ws.Range("C6").Value = 46.558181111530594
$ws.Range("D6").Value = 57.16541403083391
$ws.Range("E6").Value = 2976.0
$ws.Range("G6").Value = 78.60244750976562
$ws.Range("H6").Value = 5144.88720703125

$ws.Range("C7").Value = 30.1651314147698
$ws.Range("D7").Value = 37.03757710546576
$ws.Range("E7").Value = 2976.0
$ws.Range("G7").Value = 72.77883911132812
$ws.Range("H7").Value = 7777.89111328125

$ws.Range("C9").Value = 18.889209435992342
$ws.Range("D9").Value = 23.192690318951044
$ws.Range("E9").Value = 2976.0
$ws.Range("G9").Value = 40.355281829833984
$ws.Range("H9").Value = 4005.12353515625

$ws.Range("C10").Value = 26.7593132475371
$ws.Range("D10").Value = 32.85581992774881
$ws.Range("E10").Value = 2976.0
$ws.Range("G10").Value = 27.927446365356445
$ws.Range("H10").Value = 2628.465576171875

$ws.Range("C11").Value = 6.103590802479816
$ws.Range("D11").Value = 7.494156486244612
$ws.Range("E11").Value = 2976.0
$ws.Range("G11").Value = 32.11246109008789
$ws.Range("H11").Value = 2248.246826171875

$ws.Range("C12").Value = 21.70458971941343
$ws.Range("D12").Value = 26.649491610065585
$ws.Range("E12").Value = 2976.0
$ws.Range("G12").Value = 98.20337677001953
$ws.Range("H12").Value = 6395.8779296875

$ws.Range("C13").Value = 16.378528225806452
$ws.Range("D13").Value = 20.11000612858803
$ws.Range("E13").Value = 2976.0
$ws.Range("G13").Value = 29.159509658813477
$ws.Range("H13").Value = 3319.7470703125

$ws.Range("C14").Value = 19.696025640934064
$ws.Range("D14").Value = 24.183321226648587
$ws.Range("E14").Value = 2976.0
$ws.Range("G14").Value = 217.64988708496094
$ws.Range("H14").Value = 3285.77734375

$ws.Range("C15").Value = 5.8323734088610575
$ws.Range("D15").Value = 7.1611482968093245
$ws.Range("E15").Value = 2976.0
$ws.Range("G15").Value = 41.82110595703125
$ws.Range("H15").Value = 1861.8985595703125

$ws.Range("C16").Value = 29.662635196921645
$ws.Range("D16").Value = 36.4205983877182
$ws.Range("E16").Value = 2976.0
$ws.Range("G16").Value = 41.29185485839844
$ws.Range("H16").Value = 910.5149536132812

$ws.Range("C17").Value = 17.02524992983828
$ws.Range("D17").Value = 20.904069359584522
$ws.Range("E17").Value = 2976.0
$ws.Range("G17").Value = 34.46185302734375
$ws.Range("H17").Value = 731.6424560546875

$ws.Range("C18").Value = 9.176939534564172
$ws.Range("D18").Value = 11.267698496580124
$ws.Range("E18").Value = 2976.0
$ws.Range("G18").Value = 2.5915706157684326
$ws.Range("H18").Value = 901.4158935546875

$ws.Range("C19").Value = 10.194652644616943
$ws.Range("D19").Value = 12.517274501583268
$ws.Range("E19").Value = 2976.0
$ws.Range("G19").Value = 4.381046295166016
$ws.Range("H19").Value = 751.0364990234375

$ws.Range("C20").Value = 13.526305961352522
$ws.Range("D20").Value = 16.60797001982248
$ws.Range("E20").Value = 2976.0
$ws.Range("G20").Value = 6.643187999725342
$ws.Range("H20").Value = 996.4782104492188

$ws.Range("C21").Value = 6.873800071497118
$ws.Range("D21").Value = 8.439840540770561
$ws.Range("E21").Value = 2976.0
$ws.Range("G21").Value = 4.768509864807129
$ws.Range("H21").Value = 843.9840698242188

$ws.Range("C22").Value = 15.459389871166598
$ws.Range("D22").Value = 18.98146366303967
$ws.Range("E22").Value = 2976.0
$ws.Range("G22").Value = 25.646066665649414
$ws.Range("H22").Value = 1793.748291015625

$ws.Range("C23").Value = 6.701229082961237
$ws.Range("D23").Value = 8.22795319877645
$ws.Range("E23").Value = 2976.0
$ws.Range("G23").Value = 33.3928337097168
$ws.Range("H23").Value = 641.7803344726562

$ws.Range("C24").Value = 17.353111109425946
$ws.Range("D24").Value = 21.306626609576647
$ws.Range("E24").Value = 2976.0
$ws.Range("G24").Value = 54.225364685058594
$ws.Range("H24").Value = 1704.5301513671875

$ws.Range("C25").Value = 33.69095706651288
$ws.Range("D25").Value = 41.36668263648146
$ws.Range("E25").Value = 2976.0
$ws.Range("G25").Value = 75.42524719238281
$ws.Range("H25").Value = 827.3336791992188

$ws.Range("C26").Value = 15.005760780906165
$ws.Range("D26").Value = 18.424485367472453
$ws.Range("E26").Value = 2976.0
$ws.Range("G26").Value = 22.562908172607422
$ws.Range("H26").Value = 1289.7139892578125

$ws.Range("C27").Value = 29.4314044783673
$ws.Range("D27").Value = 36.13668708503246
$ws.Range("E27").Value = 2976.0
$ws.Range("G27").Value = 142.19786071777344
$ws.Range("H27").Value = 2710.25146484375

$ws.Range("C28").Value = 10.527313981325396
$ws.Range("D28").Value = 12.925725298623243
$ws.Range("E28").Value = 2976.0
$ws.Range("G28").Value = 0.6462862491607666
$ws.Range("H28").Value = 3102.174072265625

$ws.Range("C29").Value = 17.38517709492996
$ws.Range("D29").Value = 21.345998065045443
$ws.Range("E29").Value = 2976.0
$ws.Range("H29").Value = 800.4749145507812

# Row 8 ("Carne de pollo") values removed as outliers
$ws.Range("C8:H8").ClearContents()

